# Update "2020-08-07" refresh of the "Fonds de solidarite" volet 2 dataset:
# several (nombre_aides, montant_total) pairs were revised upward.
# Values are written with a leading apostrophe so Excel stores them as text
# (matching the source data, which keeps these columns as strings, e.g.
# region codes elsewhere in the sheet have leading zeros) and so that the
# exact decimal text (e.g. "943828.79") is preserved instead of being
# reformatted as a floating point number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = "'411"
$ws.Cells.Item(2, 4).Value = "'943828.79"
$ws.Cells.Item(14, 3).Value = "'206"
$ws.Cells.Item(14, 4).Value = "'528862.00"
$ws.Cells.Item(16, 3).Value = "'464"
$ws.Cells.Item(16, 4).Value = "'1642708.75"
$ws.Cells.Item(17, 3).Value = "'137"
$ws.Cells.Item(17, 4).Value = "'418415.33"
$ws.Cells.Item(20, 3).Value = "'161"
$ws.Cells.Item(20, 4).Value = "'412599.00"
$ws.Cells.Item(22, 3).Value = "'152"
$ws.Cells.Item(22, 4).Value = "'417137.26"
$ws.Cells.Item(28, 3).Value = "'243"
$ws.Cells.Item(28, 4).Value = "'623174.00"
$ws.Cells.Item(30, 3).Value = "'489"
$ws.Cells.Item(30, 4).Value = "'1897299.70"
$ws.Cells.Item(32, 3).Value = "'348"
$ws.Cells.Item(32, 4).Value = "'1098925.37"
$ws.Cells.Item(35, 3).Value = "'261"
$ws.Cells.Item(35, 4).Value = "'615150.71"
$ws.Cells.Item(36, 3).Value = "'189"
$ws.Cells.Item(36, 4).Value = "'548514.00"
$ws.Cells.Item(37, 3).Value = "'165"
$ws.Cells.Item(37, 4).Value = "'430174.14"
$ws.Cells.Item(38, 3).Value = "'5"
$ws.Cells.Item(38, 4).Value = "'11500.00"
$ws.Cells.Item(45, 3).Value = "'300"
$ws.Cells.Item(45, 4).Value = "'810006.74"
$ws.Cells.Item(47, 3).Value = "'544"
$ws.Cells.Item(47, 4).Value = "'1995430.99"
$ws.Cells.Item(48, 3).Value = "'359"
$ws.Cells.Item(48, 4).Value = "'1170796.16"
$ws.Cells.Item(51, 3).Value = "'3260"
$ws.Cells.Item(51, 4).Value = "'7422377.36"
$ws.Cells.Item(52, 3).Value = "'22"
$ws.Cells.Item(52, 4).Value = "'138500.00"
$ws.Cells.Item(53, 3).Value = "'3776"
$ws.Cells.Item(53, 4).Value = "'12724819.79"
$ws.Cells.Item(54, 3).Value = "'11"
$ws.Cells.Item(54, 4).Value = "'54500.00"
$ws.Cells.Item(55, 3).Value = "'3859"
$ws.Cells.Item(55, 4).Value = "'11755038.49"
$ws.Cells.Item(56, 3).Value = "'54"
$ws.Cells.Item(56, 4).Value = "'143350.00"
$ws.Cells.Item(57, 3).Value = "'81"
$ws.Cells.Item(57, 4).Value = "'291689.00"
$ws.Cells.Item(73, 3).Value = "'359"
$ws.Cells.Item(73, 4).Value = "'883635.70"
$ws.Cells.Item(75, 3).Value = "'873"
$ws.Cells.Item(75, 4).Value = "'2907919.89"
$ws.Cells.Item(76, 3).Value = "'495"
$ws.Cells.Item(76, 4).Value = "'1582502.87"
$ws.Cells.Item(78, 3).Value = "'32"
$ws.Cells.Item(78, 4).Value = "'123736.09"
$ws.Cells.Item(85, 3).Value = "'204"
$ws.Cells.Item(85, 4).Value = "'473071.00"
$ws.Cells.Item(87, 3).Value = "'473"
$ws.Cells.Item(87, 4).Value = "'1576000.50"
